$d = $word.ActiveDocument

# Move to the very end of the document content and append two new paragraphs.
$end = $d.Content
$end.Collapse(0)  # wdCollapseEnd

$end.InsertParagraphAfter()
$end.Collapse(0)
$end.MoveStart(1, 1)  # wdCharacter
$end.Text = "I want to continue mu studies in automation."

$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.MoveStart(1, 1)
$end.Text = "Hello my life!"
